$wb = $excel.ActiveWorkbook

# The change targets the "key" worksheet.
$ws = $wb.Worksheets.Item("key")

# F2 previously held "test_name" -> clear it out.
$ws.Range("F2").ClearContents()

# D4 was empty -> now holds "gender".
$ws.Range("D4").Value = "gender"
